$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters for the data columns present in the sheet (A, G and H are left untouched)
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @(
    @(1.02,1.051620220994897,1.05902303866095,1.062443487835026,1.070938647436477,1.048371382920678,1.056647133377737,1.061753675251373,1.065164820581668,1.073637150845718,1.022656987288773),
    @(1.02,1.05259577599296,1.059809775251081,1.063385249402204,1.071868032854722,1.048644079059168,1.057272717665412,1.062354629252366,1.065921080010815,1.074382714789951,1.022868941944081),
    @(1.02,1.053227387565568,1.060319204840972,1.063995772037306,1.072470254568115,1.048819585899043,1.057677241351699,1.062743185469584,1.066410940673796,1.074865358547082,1.023005909445344),
    @(1.02,1.053493003019805,1.060533453654576,1.064252707356446,1.072723630132222,1.048893142066973,1.057847237373633,1.062906461516984,1.066616999168737,1.075068312168168,1.023063446809675),
    @(1.02,1.053537606041014,1.060569431900373,1.064295863816385,1.072766184786927,1.048905479151309,1.05787577657533,1.062933871978397,1.066651604365101,1.075102391886517,1.02307310500894),
    @(1.02,1.053230936395502,1.060322067313118,1.063999204154021,1.072473639394408,1.048820569652346,1.057679513106174,1.062745367459999,1.06641369356311,1.07486807022742,1.023006678434873),
    @(1.02,1.051949838898909,1.059288845048633,1.062761523755677,1.071252561409455,1.048463737400333,1.056858608175626,1.06195683230851,1.065420295696933,1.073889072720648,1.022728655685827),
    @(1.02,1.049695190092307,1.057470972898209,1.0605893605825,1.069107410461354,1.047827731029815,1.05541002785822,1.060565064463566,1.063673757342059,1.072165639799769,1.022237371222895),
    @(1.02,1.048194023819818,1.056261008970491,1.059147236765545,1.067681784465728,1.047398899150389,1.054442983000596,1.059635745648256,1.06251212163434,1.071017879871248,1.021908946415029),
    @(1.02,1.047544469582129,1.055737558428926,1.058524216676374,1.067065549417121,1.047212070951339,1.05402393743549,1.059233001393866,1.062009779644931,1.070521183461397,1.021766525316179),
    @(1.02,1.047303266072897,1.05554319761148,1.058293015229259,1.06683681431109,1.047142503624034,1.053868239605291,1.059083353513375,1.061823286744834,1.070336733045496,1.021713592417176),
    @(1.02,1.04735500187284,1.055584885389303,1.058342598931923,1.066885871416176,1.047157433783752,1.05390163937224,1.059115455774428,1.061863285595633,1.070376296244986,1.021724948117458),
    @(1.02,1.047524530196008,1.055721491027367,1.058505101079958,1.067046638789548,1.047206323980522,1.054011068331225,1.059220632483604,1.061994362052325,1.070505935813909,1.02176215050405),
    @(1.02,1.047628991467103,1.055805667846875,1.058605252670245,1.067145714380948,1.047236424196062,1.054078485097633,1.059285428627188,1.062075135770126,1.070585816990065,1.021785067976642),
    @(1.02,1.048237142356029,1.056295758696555,1.05918861477246,1.067722704612122,1.047411274316379,1.054470787232398,1.059662467298979,1.062545474283745,1.071050850196506,1.021918394017774),
    @(1.02,1.048618743129019,1.056603306859409,1.059554925888151,1.068084922775036,1.047520647958714,1.054716785853564,1.059898882417367,1.06284068101466,1.071342631908297,1.022001969670976),
    @(1.02,1.048841368807515,1.056782739947523,1.059768726761687,1.068296301699419,1.047584333618815,1.054860242846763,1.060036746199838,1.063012933132285,1.071512851275243,1.022050697555418),
    @(1.02,1.048917285922787,1.056843929667887,1.05984165066905,1.068368393917504,1.047606030089229,1.054909152890282,1.060083748538431,1.063071677311851,1.071570896404833,1.022067309037699),
    @(1.02,1.0485777963799,1.056570305143022,1.059515609910058,1.068046049511108,1.047508924590646,1.054690395610808,1.059873520736201,1.062809001599675,1.07131132358543,1.021993004899902),
    @(1.02,1.047474606396627,1.055681262074534,1.058457242276129,1.066999292322582,1.047191931753517,1.053978845484448,1.059189661958144,1.061955760553845,1.070467758948829,1.021751196191466),
    @(1.02,1.046781391558721,1.0551227026612,1.057793054605185,1.066342092583048,1.046991636547458,1.053531201933598,1.058759399684474,1.0614198689697,1.069937635762477,1.021598980177827),
    @(1.02,1.04714883930193,1.055418765644982,1.058145034144198,1.066690397278187,1.04709791044838,1.053768531002224,1.05898751736929,1.061703900395328,1.070218639174609,1.021679689837353),
    @(1.02,1.048596298316097,1.056585217057791,1.059533374682573,1.068063614347042,1.047514222216374,1.054702320317802,1.059884980690746,1.062823315969036,1.071325470382354,1.021997055757065),
    @(1.02,1.050277733441323,1.057940597701079,1.06114986754957,1.069661200178359,1.047993006904874,1.055784757563107,1.060925133313058,1.064124804058746,1.072610982375074,1.022364540539766)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $row).Value = $rowValues[$j]
    }
}

Write-Host "vm_pu values updated for rows 2-25"
